$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "Stephon Castle"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "San Antonio Spurs"

# Row 6
$ws.Range("A6").Value = "Keegan Murray"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Sacramento Kings"

# Row 10
$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"

# Row 11
$ws.Range("A11").Value = "Daniel Gafford"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Dallas Mavericks"

# Row 14
$ws.Range("A14").Value = "Franz Wagner"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Orlando Magic"

# Row 15
$ws.Range("A15").Value = "Karl-Anthony Towns"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "New York Knicks"

# Row 16
$ws.Range("A16").Value = "Stephen Curry"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Golden State Warriors"
